$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New timestamp value applied across all data rows (column I)
$newTimestamp = 1706178604

# Column E ("reservasjon") updates per row, and all rows get the refreshed timestamp.
$updates = @{
    2  = "NEI"
    3  = "NEI"
    4  = "NEI"
    5  = "NEI"
    6  = "JA"
    7  = "NEI"
    8  = $null
    9  = "NEI"
    10 = "NEI"
    11 = "NEI"
    12 = "NEI"
    13 = $null
    14 = "NEI"
    15 = "NEI"
    16 = "JA"
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    if ($null -eq $value) {
        $ws.Cells.Item($row, 5).Value = ""
    } else {
        $ws.Cells.Item($row, 5).Value = $value
    }
    $ws.Cells.Item($row, 9).Value = $newTimestamp
}
